$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the value tested in the first scenario (delTrueTest): 4 -> 0
$ws.Range("A3").Value = 0

# Remove the now-obsolete "delValorInexistenteTest" scenario block (rows 9-11),
# including its header row, table header row and data row, and the now-unused
# shared string that went with it.
$ws.Range("A9:D11").Clear()

# Move the active selection to A3
$ws.Range("A3").Select()
